$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.793.97"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "'2.628.95"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.36%  "

$ws.Range("D5").Value = "'520.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").Value = "'144.82"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("D9").Value = "'2.638.71"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("D10").Value = "'6.29"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").Value = "'0.127"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").Value = "'3.091.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.55%  "

$ws.Range("D15").Value = "'58.828.15"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "'20.79"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").Value = "'2.633.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.21%  "

$ws.Range("D19").Value = "'345.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("D21").Value = "'10.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").Value = "'6.15"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.80%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'61.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("E25").Value = "  -0.59%  "

$ws.Range("D26").Value = "'0.164"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.19%  "

$ws.Range("D27").Value = "'0.995"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "'" + "0.0" + [char]8323 + "0799"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").Value = "'7.07"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "'6.24"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.53%  "

$ws.Range("D32").Value = "'18.86"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("E33").Value = "  +2.24%  "

$ws.Range("D34").Value = "'150.15"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").Value = "'0.978"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.50%  "

$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").Value = "'36.60"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("D39").Value = "'0.836"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.34%  "

$ws.Range("D40").Value = "'3.64"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.99%  "

$ws.Range("D41").Value = "'1.42"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'277.35"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.31%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.997"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "'0.0982"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'19.58"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.599"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").Value = "'1.991.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.43%  "

$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Value = "'4.66"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.73%  "
